# Fruta / hortaliza, semanal
# Insert two new weekly price rows into the Jengibre - Vega Monumental Concepcion sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at sheet row 9 (pushes old rows 9-15 down to 10-16) ---
$ws.Rows("9:9").Insert()

$ws.Cells.Item(9, 1).Value = 11
$ws.Cells.Item(9, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(9, 3).Value = "Bíobío"
$ws.Cells.Item(9, 4).Value = 44435
$ws.Cells.Item(9, 5).Value = 8
$ws.Cells.Item(9, 6).Value = 100114007
$ws.Cells.Item(9, 7).Value = "Jengibre"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 13000
$ws.Cells.Item(9, 12).Value = 14000
$ws.Cells.Item(9, 13).Value = 13500
$ws.Cells.Item(9, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(9, 15).Value = "Perú"
$ws.Cells.Item(9, 16).Value = 1038
$ws.Cells.Item(9, 17).Value = 13
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# --- Insert second new row at sheet row 14 (pushes old rows 13-15, now at 14-16, down to 15-17) ---
$ws.Rows("14:14").Insert()

$ws.Cells.Item(14, 1).Value = 11
$ws.Cells.Item(14, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(14, 3).Value = "Bíobío"
$ws.Cells.Item(14, 4).Value = 44433
$ws.Cells.Item(14, 5).Value = 8
$ws.Cells.Item(14, 6).Value = 100114007
$ws.Cells.Item(14, 7).Value = "Jengibre"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 13000
$ws.Cells.Item(14, 12).Value = 14000
$ws.Cells.Item(14, 13).Value = 13500
$ws.Cells.Item(14, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(14, 15).Value = "Perú"
$ws.Cells.Item(14, 16).Value = 1038
$ws.Cells.Item(14, 17).Value = 13
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# Ensure the new date cells use the same date number format as the rest of column D
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
$ws.Range("D14").NumberFormat = $ws.Range("D13").NumberFormat
